$d = $word.ActiveDocument

$d.Content.Find.Execute("địa chỉ..................................", $true, $false, $false, $false, $false, $true, 1, $false, "địa chỉ `$`{diaChi`}", 2)
